# Natmi following Dr Hou advice
# Update the LR-pair table (Hspg2 -> Col13a1) with the recomputed
# ligand-receptor communication statistics, including the new
# 'FAPs' target-cluster rows (table grows from 10 to 15 data rows).

$data = @(
    @{ "A"="ECs"; "B"="Hspg2"; "C"="Col13a1"; "D"="ECs"; "E"=3; "F"=1; "G"=197.1981153333333; "H"=591.594346; "I"=0.5013718116429122; "J"=0.5302871392468994; "K"=3; "L"=1; "M"=0.3015963333333334; "N"=0.9047890000000001; "O"=0.4737402899861982; "P"=0.5687660179545673; "Q"=59.47422852477712; "R"=535.2680567229941; "S"=0.2375200274386187; "T"=0.3016093045619781 },
    @{ "A"="ECs"; "B"="Hspg2"; "C"="Col13a1"; "D"="FAPs"; "E"=3; "F"=1; "G"=197.1981153333333; "H"=591.594346; "I"=0.5013718116429122; "J"=0.5302871392468994; "K"=1; "L"=0.3333333333333333; "M"=0.01594066666666667; "N"=0.047822; "O"=0.02503921704145383; "P"=0.03006173650500096; "Q"=3.143469423823556; "R"=28.291224814412; "S"=0.01255395761019379; "T"=0.01594135225203104 },
    @{ "A"="ECs"; "B"="Hspg2"; "C"="Col13a1"; "D"="sCs"; "E"=3; "F"=1; "G"=197.1981153333333; "H"=591.594346; "I"=0.5013718116429122; "J"=0.5302871392468994; "K"=2; "L"=1; "M"=0.319091; "N"=0.638182; "O"=0.501220492972348; "P"=0.4011722455404317; "Q"=62.92414381982867; "R"=377.544862918972; "S"=0.2512978265940997; "T"=0.2127364824328902 },
    @{ "A"="FAPs"; "B"="Hspg2"; "C"="Col13a1"; "D"="ECs"; "E"=3; "F"=1; "G"=131.273506; "H"=393.820518; "I"=0.3337599622221713; "J"=0.3530087082119477; "K"=3; "L"=1; "M"=0.3015963333333334; "N"=0.9047890000000001; "O"=0.4737402899861982; "P"=0.5687660179545673; "Q"=39.59160807341134; "R"=356.324472660702; "S"=0.158115541288914; "T"=0.2007793572729953 },
    @{ "A"="FAPs"; "B"="Hspg2"; "C"="Col13a1"; "D"="FAPs"; "E"=3; "F"=1; "G"=131.273506; "H"=393.820518; "I"=0.3337599622221713; "J"=0.3530087082119477; "K"=1; "L"=0.3333333333333333; "M"=0.01594066666666667; "N"=0.047822; "O"=0.02503921704145383; "P"=0.03006173650500096; "Q"=2.092587201310667; "R"=18.833284811796; "S"=0.00835708813382838; "T"=0.01061205477023834 },
    @{ "A"="FAPs"; "B"="Hspg2"; "C"="Col13a1"; "D"="sCs"; "E"=3; "F"=1; "G"=131.273506; "H"=393.820518; "I"=0.3337599622221713; "J"=0.3530087082119477; "K"=2; "L"=1; "M"=0.319091; "N"=0.638182; "O"=0.501220492972348; "P"=0.4011722455404317; "Q"=41.888194303046; "R"=251.329165818276; "S"=0.167287332799429; "T"=0.1416172961687141 },
    @{ "A"="M1"; "B"="Hspg2"; "C"="Col13a1"; "D"="ECs"; "E"=3; "F"=1; "G"=0.278468; "H"=0.835404; "I"=0.0007079986814710624; "J"=0.0007488306814808822; "K"=3; "L"=1; "M"=0.3015963333333334; "N"=0.9047890000000001; "O"=0.4737402899861982; "P"=0.5687660179545673; "Q"=0.08398492775066667; "R"=0.755864349756; "S"=0.000335407500669947; "T"=0.0004259094448280863 },
    @{ "A"="M1"; "B"="Hspg2"; "C"="Col13a1"; "D"="FAPs"; "E"=3; "F"=1; "G"=0.278468; "H"=0.835404; "I"=0.0007079986814710624; "J"=0.0007488306814808822; "K"=1; "L"=0.3333333333333333; "M"=0.01594066666666667; "N"=0.047822; "O"=0.02503921704145383; "P"=0.03006173650500096; "Q"=0.004438965565333334; "R"=0.03995069008800001; "S"=0.00001772773265041707; "T"=0.00002251115063353858 },
    @{ "A"="M1"; "B"="Hspg2"; "C"="Col13a1"; "D"="sCs"; "E"=3; "F"=1; "G"=0.278468; "H"=0.835404; "I"=0.0007079986814710624; "J"=0.0007488306814808822; "K"=2; "L"=1; "M"=0.319091; "N"=0.638182; "O"=0.501220492972348; "P"=0.4011722455404317; "Q"=0.088856632588; "R"=0.5331397955280001; "S"=0.0003548634481506983; "T"=0.0003004100860192572 },
    @{ "A"="M2"; "B"="Hspg2"; "C"="Col13a1"; "D"="ECs"; "E"=3; "F"=1; "G"=0.2270173333333333; "H"=0.681052; "I"=0.0005771865085793579; "J"=0.0006104742535155658; "K"=3; "L"=1; "M"=0.3015963333333334; "N"=0.9047890000000001; "O"=0.4737402899861982; "P"=0.5687660179545673; "Q"=0.06846759533644445; "R"=0.616208358028; "S"=0.0002734365039505063; "T"=0.0003472170102358354 },
    @{ "A"="M2"; "B"="Hspg2"; "C"="Col13a1"; "D"="FAPs"; "E"=3; "F"=1; "G"=0.2270173333333333; "H"=0.681052; "I"=0.0005771865085793579; "J"=0.0006104742535155658; "K"=1; "L"=0.3333333333333333; "M"=0.01594066666666667; "N"=0.047822; "O"=0.02503921704145383; "P"=0.03006173650500096; "Q"=0.003618807638222222; "R"=0.032569268744; "S"=0.0000144522982617175; "T"=0.0000183519161522721 },
    @{ "A"="M2"; "B"="Hspg2"; "C"="Col13a1"; "D"="sCs"; "E"=3; "F"=1; "G"=0.2270173333333333; "H"=0.681052; "I"=0.0005771865085793579; "J"=0.0006104742535155658; "K"=2; "L"=1; "M"=0.319091; "N"=0.638182; "O"=0.501220492972348; "P"=0.4011722455404317; "Q"=0.07243918791066667; "R"=0.434635127464; "S"=0.0002892977063671342; "T"=0.0002449053271274583 },
    @{ "A"="sCs"; "B"="Hspg2"; "C"="Col13a1"; "D"="ECs"; "E"=2; "F"=1; "G"=64.34001; "H"=128.68002; "I"=0.1635830409448661; "J"=0.1153448476061565; "K"=3; "L"=1; "M"=0.3015963333333334; "N"=0.9047890000000001; "O"=0.4737402899861982; "P"=0.5687660179545673; "Q"=19.40471110263; "R"=116.42826661578; "S"=0.07749587725404498; "T"=0.06560422966453001 },
    @{ "A"="sCs"; "B"="Hspg2"; "C"="Col13a1"; "D"="FAPs"; "E"=2; "F"=1; "G"=64.34001; "H"=128.68002; "I"=0.1635830409448661; "J"=0.1153448476061565; "K"=1; "L"=0.3333333333333333; "M"=0.01594066666666667; "N"=0.047822; "O"=0.02503921704145383; "P"=0.03006173650500096; "Q"=1.02562265274; "R"=6.153735916440001; "S"=0.004095991266519531; "T"=0.003467466415945767 },
    @{ "A"="sCs"; "B"="Hspg2"; "C"="Col13a1"; "D"="sCs"; "E"=2; "F"=1; "G"=64.34001; "H"=128.68002; "I"=0.1635830409448661; "J"=0.1153448476061565; "K"=2; "L"=1; "M"=0.319091; "N"=0.638182; "O"=0.501220492972348; "P"=0.4011722455404317; "Q"=20.53031813091; "R"=82.12127252364002; "S"=0.08199117242430158; "T"=0.04627315152568068 }
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$r = 2
foreach ($row in $data) {
    $colIdx = 1
    foreach ($col in $cols) {
        $ws.Cells.Item($r, $colIdx).Value = $row[$col]
        $colIdx = $colIdx + 1
    }
    $r = $r + 1
}

# Sheet now spans A1:T16 (1 header row + 15 data rows).
Write-Host ("Dimension should now be A1:T16; last row written = {0}" -f ($r - 1))
